$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 512.7857
$ws.Range("I8").Value = 108.25
$ws.Range("J8").Value = 2940
$ws.Range("K8").Value = 324.75
$ws.Range("L8").Value = 8820
$ws.Range("M8").Value = -185.75
$ws.Range("N8").Value = -9098
$ws.Range("H9").Value = 1508.2
$ws.Range("I9").Value = 1178.2858
$ws.Range("K9").Value = 1178.2858
$ws.Range("M9").Value = -1009.2858
$ws.Range("H33").Value = 13158875
$ws.Range("I33").Value = 16667179
$ws.Range("J33").Value = 2734.75
$ws.Range("K33").Value = 16667179
$ws.Range("L33").Value = 2734.75
$ws.Range("M33").Value = -16666950
$ws.Range("N33").Value = -3192.75
$ws.Range("H74").Value = 3062.3333
$ws.Range("I74").Value = 2134.9
$ws.Range("J74").Value = 7699.5
$ws.Range("K74").Value = 2134.9
$ws.Range("L74").Value = 7699.5
$ws.Range("M74").Value = -1198.9
$ws.Range("N74").Value = -9571.5
$ws.Range("H77").Value = 3062.3333
$ws.Range("I77").Value = 2134.9
$ws.Range("J77").Value = 7699.5
$ws.Range("K77").Value = 10674.5
$ws.Range("L77").Value = 38497.5
$ws.Range("M77").Value = -5994.5
$ws.Range("N77").Value = -47857.5
$ws.Range("H98").Value = 1467.2941
$ws.Range("I98").Value = 1276.129
$ws.Range("K98").Value = 1276.129
$ws.Range("M98").Value = 221.8710000000001
$ws.Range("H122").Value = 1467.2941
$ws.Range("I122").Value = 1276.129
$ws.Range("K122").Value = 3828.387
$ws.Range("M122").Value = -1378.387
$ws.Range("H125").Value = 7806.5
$ws.Range("J125").Value = 10000
$ws.Range("L125").Value = 90000
$ws.Range("N125").Value = -94920
$ws.Range("H135").Value = 1044.069
$ws.Range("I135").Value = 859.16
$ws.Range("K135").Value = 7732.44
$ws.Range("M135").Value = -5197.44
$ws.Range("H137").Value = 11893.588
$ws.Range("I137").Value = 3117.0588
$ws.Range("K137").Value = 9351.1764
$ws.Range("M137").Value = -6801.1764
$ws.Range("H138").Value = 2023.0834
$ws.Range("J138").Value = 3004.9
$ws.Range("L138").Value = 9014.700000000001
$ws.Range("N138").Value = -19294.7
$ws.Range("H141").Value = 5884.95
$ws.Range("I141").Value = 5668.3687
$ws.Range("K141").Value = 17005.1061
$ws.Range("M141").Value = -11825.1061

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 34231
$ws.Range("I34").Value = 20184.334
$ws.Range("J34").Value = 39498.5
$ws.Range("K34").Value = 20184.334
$ws.Range("L34").Value = 39498.5
$ws.Range("M34").Value = -19913.334
$ws.Range("N34").Value = -40040.5
$ws.Range("H61").Value = 15405.863
$ws.Range("I61").Value = 3774.8333
$ws.Range("J61").Value = 19767.5
$ws.Range("K61").Value = 3774.8333
$ws.Range("L61").Value = 19767.5
$ws.Range("M61").Value = -3562.8333
$ws.Range("N61").Value = -20191.5
$ws.Range("H136").Value = 15405.863
$ws.Range("I136").Value = 3774.8333
$ws.Range("J136").Value = 19767.5
$ws.Range("K136").Value = 11324.4999
$ws.Range("L136").Value = 59302.5
$ws.Range("M136").Value = -8774.499899999999
$ws.Range("N136").Value = -64402.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H31").Value = 225
$ws.Range("I31").Value = 225
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 225
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = 27
$ws.Range("N31").ClearContents()
$ws.Range("H134").Value = 18019.95
$ws.Range("I134").Value = 8990.857
$ws.Range("K134").Value = 26972.571
$ws.Range("M134").Value = -24437.571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 5187.25
$ws.Range("I6").Value = 4083
$ws.Range("K6").Value = 4083
$ws.Range("M6").Value = -3970
$ws.Range("H14").Value = 712.5
$ws.Range("I14").Value = 712.5
$ws.Range("K14").Value = 712.5
$ws.Range("M14").Value = -542.5
$ws.Range("H15").Value = 27486.334
$ws.Range("I15").Value = 450
$ws.Range("J15").Value = 41004.5
$ws.Range("K15").Value = 450
$ws.Range("L15").Value = 41004.5
$ws.Range("M15").Value = -280
$ws.Range("N15").Value = -41344.5
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H58").Value = 15644.703
$ws.Range("I58").Value = 12899.833
$ws.Range("K58").Value = 12899.833
$ws.Range("M58").Value = -12696.833
$ws.Range("H59").Value = 43666.332
$ws.Range("J59").Value = 43666.332
$ws.Range("L59").Value = 43666.332
$ws.Range("N59").Value = -45956.332
$ws.Range("H120").Value = 49497
$ws.Range("J120").Value = 49497
$ws.Range("L120").Value = 49497
$ws.Range("N120").Value = -56755
$ws.Range("H125").Value = 107499
$ws.Range("J125").Value = 107499
$ws.Range("L125").Value = 107499
$ws.Range("N125").Value = -112419
$ws.Range("H133").Value = 76188
$ws.Range("J133").Value = 76188
$ws.Range("L133").Value = 76188
$ws.Range("N133").Value = -81248
$ws.Range("H134").Value = 45464052
$ws.Range("I134").Value = 3553.3333
$ws.Range("K134").Value = 10659.9999
$ws.Range("M134").Value = -8124.999899999999
$ws.Range("H136").Value = 15644.703
$ws.Range("I136").Value = 12899.833
$ws.Range("K136").Value = 38699.499
$ws.Range("M136").Value = -36149.499

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 523.05554
$ws.Range("J68").Value = 623
$ws.Range("L68").Value = 1869
$ws.Range("N68").Value = -3491
$ws.Range("H71").Value = 523.05554
$ws.Range("J71").Value = 623
$ws.Range("L71").Value = 5607
$ws.Range("N71").Value = -13719
$ws.Range("H129").Value = 1518.579
$ws.Range("I129").Value = 1036.3846
$ws.Range("K129").Value = 3109.1538
$ws.Range("M129").Value = 1890.8462
$ws.Range("H131").Value = 1424
$ws.Range("J131").Value = 1497.967
$ws.Range("L131").Value = 4493.901
$ws.Range("N131").Value = -14573.901

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2429
$ws.Range("I122").Value = 1908
$ws.Range("J122").Value = 3376.2727
$ws.Range("K122").Value = 5724
$ws.Range("L122").Value = 10128.8181
$ws.Range("M122").Value = -3274
$ws.Range("N122").Value = -15028.8181

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 5223.4165
$ws.Range("I93").Value = 5573.875
$ws.Range("K93").Value = 5573.875
$ws.Range("M93").Value = -4325.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 50486.625
$ws.Range("I126").Value = 57127.57
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 171382.71
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -168912.71
$ws.Range("N126").Value = -16940
